$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update J2 message-id (uuid changed)
$ws.Range("J2").Value = @"
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:f9e4935a-5ce7-4176-aaf1-30364dd260e6"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data>
  <components xmlns="http://openconfig.net/yang/platform">
   <component>
    <name>OCH-1-1</name>
    <optical-channel xmlns="http://openconfig.net/yang/terminal-device">
     <config>
      <frequency>192000000</frequency>
     </config>
    </optical-channel>
   </component>
  </components>
 </data>
</rpc-reply>
"@

# Row 3 - G3: RPC filter for target-output-power
$ws.Range("G3").Value = @"
<get>
    <filter type="subtree">
      <components xmlns="http://openconfig.net/yang/platform">
        <component>
          <name>OCH-1-1</name>
             <optical-channel xmlns="http://openconfig.net/yang/terminal-device">
               <config>
                 <target-output-power></target-output-power>
               </config>
             </optical-channel>
        </component>
      </components>
    </filter>
</get>
"@

# Row 3 - J3: rpc-reply for target-output-power
$ws.Range("J3").Value = @"
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:55b8d346-9bdb-478d-b6c8-21e357f6ef6a"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data>
  <components xmlns="http://openconfig.net/yang/platform">
   <component>
    <name>OCH-1-1</name>
    <optical-channel xmlns="http://openconfig.net/yang/terminal-device">
     <config>
      <target-output-power>-5.70</target-output-power>
     </config>
    </optical-channel>
   </component>
  </components>
 </data>
</rpc-reply>
"@

# Row 4 - G4: RPC filter for operational-mode
$ws.Range("G4").Value = @"
<get>
    <filter type="subtree">
      <components xmlns="http://openconfig.net/yang/platform">
        <component>
          <name>OCH-1-1</name>
             <optical-channel xmlns="http://openconfig.net/yang/terminal-device">
               <config>
                 <operational-mode></operational-mode>
               </config>
             </optical-channel>
        </component>
      </components>
    </filter>
</get>
"@

# Row 4 - J4: rpc-reply for operational-mode
$ws.Range("J4").Value = @"
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:0e56178c-a7f3-4120-a002-9174cf27e804"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data>
  <components xmlns="http://openconfig.net/yang/platform">
   <component>
    <name>OCH-1-1</name>
    <optical-channel xmlns="http://openconfig.net/yang/terminal-device">
     <config>
      <operational-mode>1</operational-mode>
     </config>
    </optical-channel>
   </component>
  </components>
 </data>
</rpc-reply>
"@
